$wb = $excel.ActiveWorkbook

# "NewLoanInput" is currently the selected/active tab; it loses that status
# once "Repayment schedule" is activated below.
$wsSchedule = $wb.Worksheets.Item("Repayment schedule")

# Insert a new (empty) column before column N on the "Repayment schedule"
# sheet, shifting the former N/O/P (Late / heading / Outstanding) columns
# one place to the right (-> O/P/Q).
$wsSchedule.Columns("N:N").Insert()

# Give the freshly inserted column the same width as its neighbours
# (matches the "width=11" custom width seen in the target workbook).
$wsSchedule.Columns("N:N").ColumnWidth = 10.1666666666667

# "Repayment schedule" becomes the active/visible sheet (and therefore the
# active tab of the workbook), while "NewLoanInput" is no longer selected.
$wsSchedule.Activate() | Out-Null
$wsSchedule.Range("R8").Select() | Out-Null
